$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three paragraphs of the El Salvador notes cell (S25).
# Built with explicit newlines ([char]10) so the blank line between
# paragraphs matches the existing convention used throughout this sheet.
$nl = [char]10
$salvadorNotes = "The Government of El Salvador publishes an online dashboard that reports the cumulative number of tests performed to date (""pruebas COVID19 realizadas hasta hoy"")." + $nl + $nl + "The official dashboard only provides a snapshot of the cumulative number of tests performed as of today, making it difficult to construct a historical time series. We construct a daily time series dating back to 10 April 2020 using the figures reported in [this unofficial dashboard](https://covid-19-gis-hub-el-salvador-esri-sv.hub.arcgis.com/), supplemented by figures reported on President Nayib Bukele's [official Facebook page](https://www.facebook.com/nayibbukele). We have cross-checked a sample of unofficial figures against figures reported on the President's Facebook page." + $nl + $nl + "Note that, due to the way the data is presented by the official source, the time series may be impacted by retrospective revisions made by the source – see our [FAQ here](https://ourworldindata.org/coronavirus-testing#does-your-data-reflect-retrospective-updates-made-by-the-source)."

# --- Source URL (column D) updates ---
$ws.Range("D22").Value = "https://files.ssi.dk/Data-epidemiologisk-rapport-11092020-1lyn"
$ws.Range("D23").Value = "https://www.msp.gob.do/web/wp-content/uploads/2020/09/Boletín-Especial-174.pdf"
$ws.Range("D25").Value = "https://covid-19-gis-hub-el-salvador-esri-sv.hub.arcgis.com/"
$ws.Range("D34").Value = "https://eody.gov.gr/covid-gr-daily-report-11-09-2020/"
$ws.Range("D73").Value = "https://drive.google.com/drive/folders/1skFOAw2L0sTwfnCPtIrGBdewza7mZ6-7?usp=sharing"
$ws.Range("D96").Value = "https://ddc.moph.go.th/viralpneumonia/file/situation/situation-no252-110963.pdf"
$ws.Range("D97").Value = "https://ddc.moph.go.th/viralpneumonia/file/situation/situation-no252-110963.pdf"

# --- El Salvador notes (column S, row 25) ---
$ws.Range("S25").Value = $salvadorNotes

# --- Row-level numeric / date updates (Date=C, Number of observations=G,
#     Cumulative total=H, Cumulative total per thousand=I, Daily change=J,
#     Daily change per thousand=K, 7-day smoothed=L, 7-day smoothed per
#     thousand=M, Short-term positive rate=N, Short-term tests per case=O) ---
# Row 22
$ws.Range("C22").Value = 44084
$ws.Range("G22").Value = 215
$ws.Range("H22").Value = 2830499
$ws.Range("I22").Value = 488.674
$ws.Range("J22").Value = 6925
$ws.Range("K22").Value = 1.196
$ws.Range("L22").Value = 32848
$ws.Range("M22").Value = 5.671
$ws.Range("O22").Value = 162.844

# Row 23
$ws.Range("C23").Value = 44082
$ws.Range("G23").Value = 170
$ws.Range("H23").Value = 403755
$ws.Range("I23").Value = 37.22
$ws.Range("J23").Value = 3617
$ws.Range("K23").Value = 0.333
$ws.Range("L23").Value = 3948
$ws.Range("M23").Value = 0.364
$ws.Range("O23").Value = 5.332

# Row 26
$ws.Range("N26").Value = 0.015
$ws.Range("O26").Value = 67.324

# Row 34
$ws.Range("H34").Value = 1094354
$ws.Range("I34").Value = 104.994
$ws.Range("L34").Value = 13786
$ws.Range("M34").Value = 1.323
$ws.Range("N34").Value = 0.015
$ws.Range("O34").Value = 66.37

# Row 37
$ws.Range("C37").Value = 44084
$ws.Range("G37").Value = 197
$ws.Range("H37").Value = 95485
$ws.Range("I37").Value = 279.81
$ws.Range("J37").Value = 562
$ws.Range("K37").Value = 1.647
$ws.Range("L37").Value = 521
$ws.Range("M37").Value = 1.527
$ws.Range("O37").Value = 113.969

# Row 45
$ws.Range("C45").Value = 44085
$ws.Range("G45").Value = 146
$ws.Range("H45").Value = 5818910
$ws.Range("I45").Value = 96.241
$ws.Range("J45").Value = 61422
$ws.Range("K45").Value = 1.016
$ws.Range("L45").Value = 57743
$ws.Range("M45").Value = 0.955
$ws.Range("O45").Value = 39.365

# Row 46
$ws.Range("C46").Value = 44085
$ws.Range("G46").Value = 201
$ws.Range("H46").Value = 9653269
$ws.Range("I46").Value = 159.659
$ws.Range("J46").Value = 98880
$ws.Range("K46").Value = 1.635
$ws.Range("L46").Value = 88361
$ws.Range("M46").Value = 1.461
$ws.Range("N46").Value = 0.017
$ws.Range("O46").Value = 60.238

# Row 54
$ws.Range("C54").Value = 44084
$ws.Range("G54").Value = 197
$ws.Range("H54").Value = 404766
$ws.Range("I54").Value = 646.616
$ws.Range("J54").Value = 2085
$ws.Range("K54").Value = 3.331
$ws.Range("L54").Value = 1755
$ws.Range("M54").Value = 2.804
$ws.Range("N54").Value = 0.023
$ws.Range("O54").Value = 44.191

# Row 67
$ws.Range("C67").Value = 44083
$ws.Range("G67").Value = 199
$ws.Range("H67").Value = 838251
$ws.Range("I67").Value = 154.623
$ws.Range("J67").Value = 10823
$ws.Range("K67").Value = 1.996
$ws.Range("L67").Value = 11019
$ws.Range("M67").Value = 2.033
$ws.Range("O67").Value = 104.8

# Row 73
$ws.Range("C73").Value = 44084
$ws.Range("G73").Value = 161
$ws.Range("H73").Value = 2823879
$ws.Range("I73").Value = 25.77
$ws.Range("J73").Value = 38573
$ws.Range("K73").Value = 0.352
$ws.Range("L73").Value = 34986
$ws.Range("M73").Value = 0.319
$ws.Range("O73").Value = 13.094

# Row 76
$ws.Range("C76").Value = 44083
$ws.Range("G76").Value = 193
$ws.Range("H76").Value = 2210452
$ws.Range("I76").Value = 216.781
$ws.Range("J76").Value = 21332
$ws.Range("K76").Value = 2.092
$ws.Range("L76").Value = 16342
$ws.Range("M76").Value = 1.603
$ws.Range("O76").Value = 43.135

# Row 94
$ws.Range("C94").Value = 44084
$ws.Range("G94").Value = 229
$ws.Range("H94").Value = 1162236
$ws.Range("I94").Value = 134.291
$ws.Range("J94").Value = 11139
$ws.Range("K94").Value = 1.287
$ws.Range("L94").Value = 12203
$ws.Range("M94").Value = 1.41
$ws.Range("N94").Value = 0.03
$ws.Range("O94").Value = 33.604

# Row 96
$ws.Range("C96").Value = 44085
$ws.Range("G96").Value = 178
$ws.Range("H96").Value = 430235
$ws.Range("I96").Value = 6.164
$ws.Range("J96").Value = 1173
$ws.Range("K96").Value = 0.017
$ws.Range("L96").Value = 983
$ws.Range("O96").Value = 229.367

# Row 97
$ws.Range("C97").Value = 44085
$ws.Range("G97").Value = 83
$ws.Range("H97").Value = 857645
$ws.Range("I97").Value = 12.287
$ws.Range("J97").Value = 1173
$ws.Range("K97").Value = 0.017
$ws.Range("L97").Value = 983
$ws.Range("O97").Value = 229.367

# Row 102
$ws.Range("C102").Value = 44085
$ws.Range("G102").Value = 139
$ws.Range("H102").Value = 1806428
$ws.Range("I102").Value = 41.305
$ws.Range("J102").Value = 27945
$ws.Range("K102").Value = 0.639
$ws.Range("L102").Value = 22703
$ws.Range("M102").Value = 0.519
$ws.Range("N102").Value = 0.129
$ws.Range("O102").Value = 7.742

# Row 104
$ws.Range("C104").Value = 44084
$ws.Range("G104").Value = 163
$ws.Range("H104").Value = 16145887
$ws.Range("I104").Value = 237.838
$ws.Range("J104").Value = 205659
$ws.Range("K104").Value = 3.029
$ws.Range("L104").Value = 187241
$ws.Range("M104").Value = 2.758
$ws.Range("N104").Value = 0.013
$ws.Range("O104").Value = 79.229

# Row 109
$ws.Range("C109").Value = 44085
$ws.Range("G109").Value = 155
$ws.Range("H109").Value = 128662
$ws.Range("I109").Value = 6.999
$ws.Range("J109").Value = 1133
$ws.Range("K109").Value = 0.062
$ws.Range("L109").Value = 1099
$ws.Range("M109").Value = 0.06
$ws.Range("N109").Value = 0.09
$ws.Range("O109").Value = 11.133

